# Rename the header row (row 1) string suffixes from _old/_new to the
# respective input file format versions _FV2404/_FV2410, freeze the header
# row, and turn the data range into a native Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header labels (row 1, columns A:U) -------------------------
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headersFV2404.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $headersFV2410.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into a native Excel Table -------------------
$tableRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
